# Apply "change data base and models" edit:
#  - Sheet1 (users): remove the old "ID" count column (column G); the old
#    "admin" boolean column (H) shifts left into G.
#  - Sheet2 (assets): insert two new columns right after the username
#    column: "id" and "כמות חיפושים", with sample values for the two rows.

$wb = $excel.ActiveWorkbook

# ---- Sheet1: "משתמשים" ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G1").EntireColumn.Delete()

# ---- Sheet2: "נכסים" ----
$ws2 = $wb.Worksheets.Item(2)

# Insert two blank columns before the old column B ("עיר")
$ws2.Range("B1").EntireColumn.Insert()
$ws2.Range("B1").EntireColumn.Insert()

# New header cells
$ws2.Range("B1").Value = "id"
$ws2.Range("C1").Value = "כמות חיפושים"

# New data values for the two existing data rows
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 5

$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 2

# ---- Restore / update the selected cell on each sheet ----
# (sheet2 selected first, then sheet1 last so sheet1 ends up the active tab)
$ws2.Activate()
$ws2.Range("H6").Select()

$ws1.Activate()
$ws1.Range("K8").Select()
